# Insert a new weekly data row at row 124 (shifts existing rows 124..202 down to 125..203),
# then populate the new row 124 with the new week's price data for
# Acelga - Primera - Vega Monumental Concepción.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 124; Excel copies formatting
# (incl. the date number format used by column D) from the row context.
$ws.Rows.Item(124).Insert()

$ws.Cells.Item(124, 1).Value = 11
$ws.Cells.Item(124, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(124, 3).Value = "Bíobío"
$ws.Cells.Item(124, 4).Value = 44572
$ws.Cells.Item(124, 5).Value = 8
$ws.Cells.Item(124, 6).Value = 100112009
$ws.Cells.Item(124, 7).Value = "Acelga"
$ws.Cells.Item(124, 8).Value = "Sin especificar"
$ws.Cells.Item(124, 9).Value = "Primera"
$ws.Cells.Item(124, 10).Value = 280
$ws.Cells.Item(124, 11).Value = 600
$ws.Cells.Item(124, 12).Value = 650
$ws.Cells.Item(124, 13).Value = 627
$ws.Cells.Item(124, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(124, 15).Value = "Región de Ñuble"
$ws.Cells.Item(124, 16).Value = 627
$ws.Cells.Item(124, 17).Value = 1
$ws.Cells.Item(124, 18).Value = "Hortaliza"
